# Updates numeric results across the Weights.xlsx output sheets
# following a refresh of the aerodynamic/weight estimation calculations.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C6").Value = 63803.75951503964
$ws.Range("C7").Value = 57423.38356353567
$ws.Range("C9").Value = 11804.618883255589
$ws.Range("C10").Value = 55639.14063178405
$ws.Range("C11").Value = 51999.14063178406
$ws.Range("C12").Value = 17640.0
$ws.Range("C13").Value = 14000.0
$ws.Range("C14").Value = 37999.14063178406
$ws.Range("C15").Value = 37221.01383961456
$ws.Range("C16").Value = 319.0395011694843
$ws.Range("C18").Value = 2095.8
$ws.Range("C19").Value = 35444.25334078406
$ws.Range("C20").Value = 20734.557685951215

$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C2").Value = 6380.790023389685
$ws.Range("C3").Value = 8157.5
$ws.Range("D3").Value = 27.84467080247952
$ws.Range("C5").Value = 8157.499999999999
$ws.Range("D8").Value = 229.5986848479249
$ws.Range("D9").Value = 0.23837137023096489
$ws.Range("D10").Value = 69.28938204209426
$ws.Range("C11").Value = 10244.0
$ws.Range("D11").Value = 60.54438341410976
$ws.Range("C12").Value = 7149.0
$ws.Range("D12").Value = 12.03941790584446
$ws.Range("C13").Value = 14388.0
$ws.Range("D13").Value = 125.48931946136383
$ws.Range("C14").Value = 6463.0
$ws.Range("D14").Value = 1.2883980872111829
$ws.Range("D15").Value = 23.668071995386576

$ws = $wb.Worksheets.Item("WING")
$ws.Range("C2").Value = 6763.637424793066
$ws.Range("C3").Value = 7181.25
$ws.Range("D3").Value = 6.174378503438334
$ws.Range("C5").Value = 7181.249999999999
$ws.Range("C8").Value = 6631.0
$ws.Range("D8").Value = -1.9610368868512102
$ws.Range("C9").Value = 7561.0
$ws.Range("D9").Value = 11.78896095589172
$ws.Range("C10").Value = 8395.0
$ws.Range("D10").Value = 24.119604182609574
$ws.Range("C11").Value = 6138.0
$ws.Range("D11").Value = -9.250014237896657

$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C2").Value = 733.7908526898138
$ws.Range("C3").Value = 796.625
$ws.Range("D3").Value = 8.562950475582888
$ws.Range("C5").Value = 796.6249999999999
$ws.Range("D8").Value = 92.83423809565326
$ws.Range("D9").Value = -92.91351209825162
$ws.Range("D10").Value = -4.604970553387074
$ws.Range("D11").Value = 41.72975803496778
$ws.Range("D12").Value = 0.4373381459338955
$ws.Range("D13").Value = 107.55232835313069
$ws.Range("C14").Value = 507.0
$ws.Range("D14").Value = -30.906742957953206
$ws.Range("C15").Value = 399.0
$ws.Range("D15").Value = -45.624833215430634

$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("C2").Value = 733.7908526898138
$ws.Range("C3").Value = 673.0
$ws.Range("D3").Value = -8.284493117756428
$ws.Range("C5").Value = 672.9999999999999
$ws.Range("D8").Value = 56.039012451959714
$ws.Range("D9").Value = -31.588136025429016
$ws.Range("D10").Value = -83.10145192659999
$ws.Range("D11").Value = 2.0726815078758314
$ws.Range("C12").Value = 488.0
$ws.Range("D12").Value = -33.49603661436127
$ws.Range("D13").Value = 107.55232835313069
$ws.Range("C14").Value = 180.0
$ws.Range("D14").Value = -75.46984957087096

$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("C2").Value = 1212.35010444404
$ws.Range("D3").Value = 14.708338917032469
$ws.Range("D10").Value = 16.303037780212705
$ws.Range("D11").Value = 14.488380453145556
$ws.Range("D12").Value = 13.333598517739189
$ws.Range("D17").Value = 16.303037780212705
$ws.Range("D18").Value = 14.488380453145556
$ws.Range("D19").Value = 13.333598517739189

$ws = $wb.Worksheets.Item("POWER PLANT")
$ws.Range("C2").Value = 5296.055719413439
$ws.Range("D3").Value = 21.81392483627129
$ws.Range("D11").Value = 11.554717567328343
$ws.Range("D12").Value = 23.299306993001704
$ws.Range("D13").Value = 30.58774994848389
$ws.Range("D18").Value = 11.554717567328343
$ws.Range("D19").Value = 23.299306993001704
$ws.Range("D20").Value = 30.58774994848389

$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C2").Value = 2616.123909589771
$ws.Range("C3").Value = 2535.5160192845633
$ws.Range("D3").Value = -3.081195428463017
$ws.Range("C5").Value = 2535.516019284563
$ws.Range("C9").Value = 2535.5160192845633
$ws.Range("D9").Value = -3.0811954284630008
$ws.Range("C11").Value = 391.42357705513115
$ws.Range("C13").Value = 2144.092442229433

$ws = $wb.Worksheets.Item("SYSTEMS")
$ws.Range("C2").Value = 8677.874431809973
$ws.Range("C3").Value = 8258.362321499524
$ws.Range("D3").Value = -4.834272650600571
$ws.Range("C4").Value = 8258.362321499524
$ws.Range("C8").Value = 8258.362321499524
$ws.Range("D8").Value = -4.834272650600542
$ws.Range("C11").Value = 336.84274966573867
$ws.Range("C13").Value = 336.8427496657386
$ws.Range("C21").Value = 1033.7781352913976
$ws.Range("C23").Value = 1033.7781352913973
$ws.Range("C26").Value = 531.6934489063985
$ws.Range("C28").Value = 531.6934489063984
$ws.Range("C36").Value = 785.1150819159232
$ws.Range("C38").Value = 785.115081915923
$ws.Range("C41").Value = 3321.868948171935
$ws.Range("C43").Value = 3321.8689481719343

